# Updates the cryptos price list with the latest scraped values.
# For rows 35, 36, 38, 39 the coin identity (name/link) also moved between
# rows, so Coin (B), Link (C), Price (D) and Volume(1h) (E) are all updated
# for those rows; for the rest only Price (D) and/or Volume(1h) (E) change.
#
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the workbook's original inlineStr / text formatting), instead of
# re-interpreting strings such as "29.343.08" or "1.000" as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.343.08"
$ws.Cells.Item(2, 5).Value = "'  +0.17%  "
$ws.Cells.Item(3, 4).Value = "'1.877.54"
$ws.Cells.Item(3, 5).Value = "'  +0.21%  "
$ws.Cells.Item(4, 5).Value = "'  +0.04%  "
$ws.Cells.Item(5, 4).Value = "'0.7136"
$ws.Cells.Item(5, 5).Value = "'  +0.19%  "
$ws.Cells.Item(6, 4).Value = "'242.15"
$ws.Cells.Item(6, 5).Value = "'  -0.24%  "
$ws.Cells.Item(7, 4).Value = "'1.000"
$ws.Cells.Item(7, 5).Value = "'  +0.05%  "
$ws.Cells.Item(8, 4).Value = "'0.08079"
$ws.Cells.Item(8, 5).Value = "'  +4.44%  "
$ws.Cells.Item(9, 5).Value = "'  +0.51%  "
$ws.Cells.Item(10, 4).Value = "'25.22"
$ws.Cells.Item(10, 5).Value = "'  +0.53%  "
$ws.Cells.Item(11, 4).Value = "'0.08351"
$ws.Cells.Item(11, 5).Value = "'  -1.49%  "
$ws.Cells.Item(12, 4).Value = "'1.874.05"
$ws.Cells.Item(12, 5).Value = "'  +0.57%  "
$ws.Cells.Item(13, 5).Value = "'  +0.83%  "
$ws.Cells.Item(14, 4).Value = "'0.7189"
$ws.Cells.Item(14, 5).Value = "'  +1.02%  "
$ws.Cells.Item(15, 4).Value = "'91.45"
$ws.Cells.Item(15, 5).Value = "'  +0.15%  "
$ws.Cells.Item(16, 4).Value = "'6.257"
$ws.Cells.Item(16, 5).Value = "'  +4.56%  "
$ws.Cells.Item(17, 4).Value = "'0.000008393"
$ws.Cells.Item(17, 5).Value = "'  +0.92%  "
$ws.Cells.Item(18, 4).Value = "'29.347.37"
$ws.Cells.Item(18, 5).Value = "'  +0.17%  "
$ws.Cells.Item(19, 4).Value = "'240.90"
$ws.Cells.Item(19, 5).Value = "'  -0.73%  "
$ws.Cells.Item(20, 5).Value = "'  +0.22%  "
$ws.Cells.Item(21, 4).Value = "'2.128.22"
$ws.Cells.Item(21, 5).Value = "'  +0.13%  "
$ws.Cells.Item(22, 4).Value = "'0.9991"
$ws.Cells.Item(22, 5).Value = "'  -0.05%  "
$ws.Cells.Item(23, 4).Value = "'7.810"
$ws.Cells.Item(23, 5).Value = "'  +0.22%  "
$ws.Cells.Item(24, 5).Value = "'  +0.04%  "
$ws.Cells.Item(25, 5).Value = "'  -2.33%  "
$ws.Cells.Item(26, 4).Value = "'163.21"
$ws.Cells.Item(26, 5).Value = "'  +0.24%  "
$ws.Cells.Item(27, 5).Value = "'  +0.61%  "
$ws.Cells.Item(28, 5).Value = "'  +0.23%  "
$ws.Cells.Item(29, 4).Value = "'1.505"
$ws.Cells.Item(29, 5).Value = "'  -0.25%  "
$ws.Cells.Item(30, 5).Value = "'  +0.08%  "
$ws.Cells.Item(31, 4).Value = "'4.337"
$ws.Cells.Item(31, 5).Value = "'  +0.17%  "
$ws.Cells.Item(32, 4).Value = "'1.202"
$ws.Cells.Item(32, 5).Value = "'  -5.88%  "
$ws.Cells.Item(33, 4).Value = "'0.05374"
$ws.Cells.Item(33, 5).Value = "'  +2.21%  "
$ws.Cells.Item(34, 4).Value = "'1.950"
$ws.Cells.Item(34, 5).Value = "'  +1.52%  "
$ws.Cells.Item(35, 2).Value = "'ImmutableX"
$ws.Cells.Item(35, 3).Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).Value = "'0.7514"
$ws.Cells.Item(35, 5).Value = "'  +0.99%  "
$ws.Cells.Item(36, 2).Value = "'ARBITRUM"
$ws.Cells.Item(36, 3).Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "'1.178"
$ws.Cells.Item(36, 5).Value = "'  +0.44%  "
$ws.Cells.Item(37, 4).Value = "'2.700"
$ws.Cells.Item(37, 5).Value = "'  +0.70%  "
$ws.Cells.Item(38, 2).Value = "'VeChain"
$ws.Cells.Item(38, 3).Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.01880"
$ws.Cells.Item(38, 5).Value = "'  +1.14%  "
$ws.Cells.Item(39, 2).Value = "'Maker"
$ws.Cells.Item(39, 3).Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(39, 4).Value = "'1.282.83"
$ws.Cells.Item(39, 5).Value = "'  +10.23%  "
$ws.Cells.Item(40, 4).Value = "'2.737"
$ws.Cells.Item(40, 5).Value = "'  +0.68%  "
$ws.Cells.Item(41, 4).Value = "'6.595"
$ws.Cells.Item(41, 5).Value = "'  +3.90%  "
$ws.Cells.Item(42, 4).Value = "'110.40"
$ws.Cells.Item(42, 5).Value = "'  +3.58%  "
$ws.Cells.Item(43, 4).Value = "'0.8923"
$ws.Cells.Item(43, 5).Value = "'  +0.34%  "
$ws.Cells.Item(44, 4).Value = "'73.21"
$ws.Cells.Item(44, 5).Value = "'  +0.44%  "
$ws.Cells.Item(45, 5).Value = "'  +8.69%  "
$ws.Cells.Item(46, 4).Value = "'1.0000"
$ws.Cells.Item(46, 5).Value = "'  +0.03%  "
$ws.Cells.Item(47, 4).Value = "'2.021.64"
$ws.Cells.Item(47, 5).Value = "'  +0.04%  "
$ws.Cells.Item(48, 4).Value = "'1.801"
$ws.Cells.Item(48, 5).Value = "'  -0.16%  "
$ws.Cells.Item(49, 5).Value = "'  +0.15%  "
$ws.Cells.Item(50, 4).Value = "'9.479"
$ws.Cells.Item(50, 5).Value = "'  +1.07%  "
$ws.Cells.Item(51, 4).Value = "'0.4365"
$ws.Cells.Item(51, 5).Value = "'  +1.49%  "
